$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.144.46'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '2.433.50'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  -0.22%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '574.52'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.86%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '140.71'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').Value = '2.420.29'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('E11').Value = '  +1.43%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '5.14'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.78%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.340'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.06%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '26.25'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').Value = '2.883.53'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '61.155.23'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').Value = '2.429.66'
$ws.Range('E18').Value = '  -0.24%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.63'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.21%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.27'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +3.00%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '324.26'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.44%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.06'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.90%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.11'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('E24').Value = '  -0.13%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '1.92'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.33%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '65.10'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.81%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.96'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -3.49%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '573.88'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -6.56%  '
$ws.Range('D29').Value = '2.568.47'
$ws.Range('E29').Value = '  +0.44%  '
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').Value = '0.0₃0917'
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('E32').Value = '  -0.87%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.35'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -4.53%  '
$ws.Range('E34').Value = '  -0.89%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.132'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -7.24%  '
$ws.Range('E36').Value = '  +0.13%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.64'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -5.09%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.370'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '151.45'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('E40').Value = '  -2.50%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '18.33'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.36%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.13'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('E43').Value = '  +0.01%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '41.68'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('E45').Value = '  -4.73%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.36'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -3.45%  '
$ws.Range('D47').Value = '0.0₆0291'
$ws.Range('E47').Value = '  +26.08%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '141.68'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.43%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '3.54'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.69%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.594'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.35%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0508'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.51%  '
